$d = $word.ActiveDocument

# Replace the full greeting text (spread across 3 runs) with "Coming soon"
$range = $d.Content
$range.Find.Execute("Hola  que tal como estas", $false, $false, $false, $false, $false, $true, 1, $false, "Coming soon", 2)
